# Applies the cryptos price/volume update described in the commit
# "Updated cryptos list on Fri Jul 28 02:52:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.198.79"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.860.88"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7137"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.13"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3082"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07689"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.00"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08321"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.917.40"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7162"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.75"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "29.312.27"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.961"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "2.172.48"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.69"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007799"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.987"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1614"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.76"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.898"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.58"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.341"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.440"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.494"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.247"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05184"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7903"
$ws.Range("E34").Value = "  +9.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.924"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.691"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "1.176.63"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.240"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9023"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.82"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "2.069.65"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.40"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5201"
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.774"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.341"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.010"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.061"
$ws.Range("E51").Value = "  +0.48%  "
